$d = $word.ActiveDocument

# The document previously had every bullet point struck through
# (strikethrough direct character formatting applied both to the run
# text and to the paragraph mark itself, i.e. w:pPr/w:rPr and each
# w:r/w:rPr). This commit removes that strikethrough formatting so the
# requirements read normally again.
#
# Iterating Paragraphs (rather than Document.Content) and using each
# paragraph's own Range ensures the paragraph-mark formatting is
# cleared too, not just the visible run text.
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.StrikeThrough = $false
}

Write-Host "Removed strikethrough formatting from $($d.Paragraphs.Count) paragraphs."
